$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.543685793876648
$ws.Range("B1").Value = 1.666494488716125
$ws.Range("C1").Value = 1.865381121635437
$ws.Range("D1").Value = 2.887477397918701
$ws.Range("E1").Value = 3.854363441467285
